$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Sheet, $CellRef, $Text) {
    $rng = $Sheet.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

Set-TextValue $ws "D2" "26.805.18"
Set-TextValue $ws "E2" "  -3.14%  "
Set-TextValue $ws "D3" "1.856.17"
Set-TextValue $ws "E3" "  -2.13%  "
Set-TextValue $ws "D4" "1.001"
Set-TextValue $ws "E4" "  +0.11%  "
Set-TextValue $ws "D5" "305.10"
Set-TextValue $ws "D6" "1.000"
Set-TextValue $ws "E6" "  +0.09%  "
Set-TextValue $ws "D7" "0.5083"
Set-TextValue $ws "E7" "  -3.23%  "
Set-TextValue $ws "E8" "  -3.74%  "
Set-TextValue $ws "D9" "0.07126"
Set-TextValue $ws "E9" "  -1.58%  "
Set-TextValue $ws "D10" "20.73"
Set-TextValue $ws "E10" "  -1.72%  "
Set-TextValue $ws "D11" "0.8866"
Set-TextValue $ws "E11" "  -1.72%  "
Set-TextValue $ws "D12" "1.858.18"
Set-TextValue $ws "E12" "  -2.08%  "
Set-TextValue $ws "D13" "0.07492"
Set-TextValue $ws "E13" "  -1.87%  "
Set-TextValue $ws "D14" "5.236"
Set-TextValue $ws "E14" "  -3.64%  "
Set-TextValue $ws "D15" "90.33"
Set-TextValue $ws "E15" "  -1.55%  "
Set-TextValue $ws "D16" "1.002"
Set-TextValue $ws "E16" "  +0.23%  "
Set-TextValue $ws "D17" "0.000008527"
Set-TextValue $ws "E17" "  -1.57%  "
Set-TextValue $ws "D18" "14.02"
Set-TextValue $ws "E18" "  -2.04%  "
Set-TextValue $ws "D19" "0.9999"
Set-TextValue $ws "E19" "  +0.05%  "
Set-TextValue $ws "D20" "26.853.23"
Set-TextValue $ws "E20" "  -3.06%  "
Set-TextValue $ws "D21" "4.996"
Set-TextValue $ws "E21" "  -2.90%  "
Set-TextValue $ws "D22" "2.103.87"
Set-TextValue $ws "E22" "  -1.53%  "
Set-TextValue $ws "D23" "10.25"
Set-TextValue $ws "E23" "  -5.17%  "
Set-TextValue $ws "D24" "6.437"
Set-TextValue $ws "E24" "  -2.42%  "
Set-TextValue $ws "D25" "1.818"
Set-TextValue $ws "E25" "  -2.11%  "
Set-TextValue $ws "D26" "145.90"
Set-TextValue $ws "E26" "  -4.71%  "
Set-TextValue $ws "D27" "17.81"
Set-TextValue $ws "E27" "  -2.39%  "
Set-TextValue $ws "D28" "2.041"
Set-TextValue $ws "E28" "  -6.25%  "
Set-TextValue $ws "D29" "112.76"
Set-TextValue $ws "E29" "  -1.15%  "
Set-TextValue $ws "D30" "4.616"
Set-TextValue $ws "E30" "  -4.60%  "
Set-TextValue $ws "D31" "4.653"
Set-TextValue $ws "E31" "  -3.46%  "
Set-TextValue $ws "D32" "0.09213"
Set-TextValue $ws "E32" "  +0.65%  "
Set-TextValue $ws "E33" "  -3.01%  "
Set-TextValue $ws "D34" "3.067"
Set-TextValue $ws "E34" "  -2.86%  "
Set-TextValue $ws "D35" "1.147"
Set-TextValue $ws "E35" "  -6.16%  "
Set-TextValue $ws "D36" "0.7310"
Set-TextValue $ws "E36" "  -5.45%  "
Set-TextValue $ws "D37" "3.193"
Set-TextValue $ws "E37" "  +3.85%  "
Set-TextValue $ws "D38" "0.02012"
Set-TextValue $ws "E38" "  -3.34%  "
Set-TextValue $ws "D39" "2.455"
Set-TextValue $ws "E39" "  -4.46%  "
Set-TextValue $ws "D40" "1.072"
Set-TextValue $ws "E40" "  -1.59%  "
Set-TextValue $ws "D41" "0.5285"
Set-TextValue $ws "E41" "  -5.00%  "
Set-TextValue $ws "D42" "117.73"
Set-TextValue $ws "E42" "  +0.23%  "
Set-TextValue $ws "D43" "6.452"
Set-TextValue $ws "E43" "  -3.69%  "
Set-TextValue $ws "D44" "8.338"
Set-TextValue $ws "E44" "  -4.38%  "
Set-TextValue $ws "D45" "0.1471"
Set-TextValue $ws "E45" "  -2.56%  "
Set-TextValue $ws "D46" "1.000"
Set-TextValue $ws "E46" "  +0.12%  "
Set-TextValue $ws "D47" "0.4633"
Set-TextValue $ws "E47" "  -3.52%  "
Set-TextValue $ws "D48" "9.915"
Set-TextValue $ws "E48" "  -4.28%  "
Set-TextValue $ws "D49" "1.552"
Set-TextValue $ws "E49" "  -2.61%  "
Set-TextValue $ws "D50" "36.94"
Set-TextValue $ws "E50" "  -0.26%  "
Set-TextValue $ws "E51" "  -4.71%  "
